$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the two comparison-year columns in the header row
$ws.Range("C3").Value = "CURRENT YEAR"
$ws.Range("D3").Value = "PREVIOUS YEAR"

# Clear out last year's figures -> replace with placeholder 1s; the
# "% Difference" column (E) keeps its formulas and recalculates itself
$ws.Range("C4:D10").Value = 1

# Remove the helper column (H) that held the plain-text summary notes
$ws.Range("H:H").Delete()

# Clear the leftover capex-percentage scratch calculation, keeping formatting
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()

# Match the saved selection/view state
$ws.Range("D3").Select() | Out-Null
